$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at position 11, shifting existing rows 11-25 down to 12-26
$ws.Rows.Item(11).Insert()

# Populate the newly inserted row 11 with the new weekly record
$ws.Cells.Item(11, 1).Value = 7
$ws.Cells.Item(11, 2).Value = "Terminal Hortofrutícola Agro Chillán"
$ws.Cells.Item(11, 3).Value = "Ñuble"
$ws.Cells.Item(11, 4).Value = "2022-08-24"
$ws.Cells.Item(11, 5).Value = 16
$ws.Cells.Item(11, 6).Value = 100112043
$ws.Cells.Item(11, 7).Value = "Pepino dulce"
$ws.Cells.Item(11, 8).Value = "Cultivar IV Región"
$ws.Cells.Item(11, 9).Value = "Primera"
$ws.Cells.Item(11, 10).Value = 80
$ws.Cells.Item(11, 11).Value = 16000
$ws.Cells.Item(11, 12).Value = 17000
$ws.Cells.Item(11, 13).Value = 16500
$ws.Cells.Item(11, 14).Value = "$/bandeja 18 kilos"
$ws.Cells.Item(11, 15).Value = "Provincia de Limarí"
$ws.Cells.Item(11, 16).Value = 917
$ws.Cells.Item(11, 17).Value = 18
$ws.Cells.Item(11, 18).Value = "Hortaliza"
